{"js": "// Update the worksheet date and every division-problem answer cell.\n//\n// The document is a single title paragraph (\"2023-11-09 Thursday\") followed\n// by one 20x5 table where only every 4th row (0, 4, 8, 12, 16) actually holds\n// text; the rest are blank spacer rows. We replace the title text and then\n// walk the table's non-empty cells in row-major (document) order, rewriting\n// each one positionally \u2014 NOT via global text search-and-replace, because a\n// couple of the new answers are identical to OLD answers that live elsewhere\n// in the table (e.g. \"31\u00f77=4, 3\" is both an old value at (row 8, col 0) and\n// the new value for (row 12, col 4)), which would make a naive find/replace\n// clobber the wrong cell or double-fire.\n\nconst titleOld = \"2023-11-09 Thursday\";\nconst titleNew = \"2023-11-10 Friday\";\n\n// Ordered (old -> new) answers, in the same left-to-right / top-to-bottom\n// order the cells appear in the document.\nconst newAnswers = [\n  \"31\u00f79=3, 4\",\n  \"44\u00f75=8, 4\",\n  \"40\u00f73=13, 1\",\n  \"35\u00f74=8, 3\",\n  \"13\u00f78=1, 5\",\n  \"18\u00f74=4, 2\",\n  \"47\u00f79=5, 2\",\n  \"15\u00f78=1, 7\",\n  \"17\u00f79=1, 8\",\n  \"72\u00f74=18, 0\",\n  \"31\u00f77=4, 3\",\n  \"99\u00f79=11, 0\",\n  \"94\u00f79=10, 4\",\n  \"80\u00f76=13, 2\",\n  \"64\u00f78=8, 0\",\n  \"65\u00f77=9, 2\",\n  \"19\u00f76=3, 1\",\n  \"47\u00f73=15, 2\",\n  \"48\u00f74=12, 0\",\n  \"25\u00f79=2, 7\",\n  \"20\u00f78=2, 4\",\n  \"15\u00f72=7, 1\",\n  \"44\u00f74=11, 0\",\n  \"12\u00f72=6, 0\",\n  \"33\u00f75=6, 3\",\n];\n\n// --- Title paragraph -------------------------------------------------\nconst firstPara = context.document.body.paragraphs.getFirst();\nfirstPara.load(\"text\");\nawait context.sync();\n\nif (firstPara.text.trim() === titleOld) {\n  firstPara.insertText(titleNew, \"Replace\");\n} else {\n  // Fall back to a scoped search within just this paragraph if the text\n  // doesn't match exactly (e.g. surrounding whitespace differences).\n  const hits = firstPara.search(titleOld, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(titleNew, \"Replace\");\n  }\n}\nawait context.sync();\n\n// --- Table of answers --------------------------------------------------\nconst table = context.document.body.tables.getFirst();\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nlet answerIdx = 0;\nfor (let r = 0; r < grid.length; r++) {\n  const row = grid[r];\n  for (let c = 0; c < row.length; c++) {\n    const cellText = (row[c] || \"\").trim();\n    if (cellText.length === 0) continue;\n    if (answerIdx < newAnswers.length) {\n      table.getCell(r, c).value = newAnswers[answerIdx];\n      answerIdx++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every division-problem answer cell.\n#\n# The document is a single title paragraph (\"2023-11-09 Thursday\") followed\n# by one 20x5 table where only every 4th row (rows 1, 5, 9, 13, 17 in COM's\n# 1-based indexing) actually holds text; the rest are blank spacer rows. We\n# replace the title text and then walk the table's non-empty cells in\n# row-major (document) order, rewriting each one positionally by\n# row/column index - NOT via a global Find/Replace - because a couple of the\n# new answers are identical to OLD answers that live elsewhere in the table\n# (e.g. \"31div7=4, 3\" is both the old value at row 9 col 1 and the new value\n# for row 13 col 5). A sequential text-based replace-all would clobber the\n# wrong cell or fire twice, so addressing by cell position avoids that\n# entirely.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph ----------------------------------------------------\n$titleOld = \"2023-11-09 Thursday\"\n$titleNew = \"2023-11-10 Friday\"\n\n$titlePara = $d.Paragraphs.First\n$titleRange = $titlePara.Range\n$titleText = $titleRange.Text.TrimEnd([char]13, [char]7)\n\nif ($titleText -eq $titleOld) {\n    $titleRange.Text = $titleNew\n} else {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $titleOld\n    $find.Replacement.Text = $titleNew\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# --- Table of answers, in document order --------------------------------\n$newAnswers = @(\n    \"31\u00f79=3, 4\",\n    \"44\u00f75=8, 4\",\n    \"40\u00f73=13, 1\",\n    \"35\u00f74=8, 3\",\n    \"13\u00f78=1, 5\",\n    \"18\u00f74=4, 2\",\n    \"47\u00f79=5, 2\",\n    \"15\u00f78=1, 7\",\n    \"17\u00f79=1, 8\",\n    \"72\u00f74=18, 0\",\n    \"31\u00f77=4, 3\",\n    \"99\u00f79=11, 0\",\n    \"94\u00f79=10, 4\",\n    \"80\u00f76=13, 2\",\n    \"64\u00f78=8, 0\",\n    \"65\u00f77=9, 2\",\n    \"19\u00f76=3, 1\",\n    \"47\u00f73=15, 2\",\n    \"48\u00f74=12, 0\",\n    \"25\u00f79=2, 7\",\n    \"20\u00f78=2, 4\",\n    \"15\u00f72=7, 1\",\n    \"44\u00f74=11, 0\",\n    \"12\u00f72=6, 0\",\n    \"33\u00f75=6, 3\"\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$answerIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cellRange = $cell.Range\n        $cellText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($cellText.Length -gt 0) {\n            if ($answerIndex -lt $newAnswers.Count) {\n                $cellRange.Text = $newAnswers[$answerIndex]\n            }\n            $answerIndex = $answerIndex + 1\n        }\n    }\n}\n"}
